$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two data columns (A and B) for rows 1-4 with the new values.
# Row 5 is unchanged by this edit.
$ws.Range("A1").Value = -0.057239506531692755
$ws.Range("B1").Value = 0.05723950650620295

$ws.Range("A2").Value = 0.0092863688319927608
$ws.Range("B2").Value = -0.0092863689010319727

$ws.Range("A3").Value = 0.055533955935311602
$ws.Range("B3").Value = -0.055533955987749163

$ws.Range("A4").Value = -0.026260855955375095
$ws.Range("B4").Value = 0.026260855862889371

# Column widths also shifted slightly in the source edit.
# (Excel COM snaps ColumnWidth to whole-pixel increments, so these are the
# closest attainable character widths to the target 14.7109375 / 15.42578125.)
$ws.Columns.Item(1).ColumnWidth = 13.833333333333334
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666
